$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N19").Value = 754857.45
$ws.Range("K27").Value = 12813.93
$ws.Range("M28").Value = 786519.07
$ws.Range("N28").Value = 336723.21
$ws.Range("O28").Value = 313936.69
